# Application.xlsx - "data types arxml generated"
#
# The DataTypes sheet gains a new "Bit Size" column in its first
# ("Base Types") table, a new Uint32 row, and its separator rows are
# consolidated. Everything below that table shifts up by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$xlShiftUp = -4162
$xlShiftDown = -4121
$xlPasteFormats = -4122
$xlLeft = -4131

# ---------------------------------------------------------------------
# Step 1: consolidate the "Base Types" table.
# Originally rows 4/6/8 hold Boolean/Uint8/Uint16 separated by blank rows
# 5 and 7. Remove those blank separators (only within columns B:D so the
# other tables' column E stays untouched), then insert one fresh blank
# row back in to hold a new "Uint32" entry.
# ---------------------------------------------------------------------
$ws.Range("B5:D5").Delete($xlShiftUp)
$ws.Range("B6:D6").Delete($xlShiftUp)
$ws.Range("B7:D7").Insert($xlShiftDown)

# ---------------------------------------------------------------------
# Step 2: make room for the new "Bit Size" column by moving the
# existing "Native Declaration" column (D) of this first table into a
# new column E. (Only this table's header/rows are touched - the
# tables further down already use both D and E.)
# ---------------------------------------------------------------------
foreach ($r in 2,3,4,5,6,8) {
    $ws.Range("D$r").Copy()
    $ws.Range("E$r").PasteSpecial($xlPasteFormats)
    $ws.Range("E$r").Value2 = $ws.Range("D$r").Value2
}
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# Step 3: fill in the new "Bit Size" column (D) for the header and the
# Boolean/Uint8/Uint16/Uint32 rows.
# ---------------------------------------------------------------------
# Header cell D2 takes the same style as C2 ("Type").
$ws.Range("C2").Copy()
$ws.Range("D2").PasteSpecial($xlPasteFormats)
$ws.Range("D2").Value2 = "Bit Size"

# Blank row D3 takes the same style as C3.
$ws.Range("C3").Copy()
$ws.Range("D3").PasteSpecial($xlPasteFormats)

# Data rows: base the new Bit Size cells on the existing "Type" column
# styling (font/border) and then left-align them.
foreach ($r in 4,5,6,7) {
    $ws.Range("C$r").Copy()
    $ws.Range("D$r").PasteSpecial($xlPasteFormats)
    $ws.Range("D$r").HorizontalAlignment = $xlLeft
}
$ws.Range("D4").Value2 = 8
$ws.Range("D5").Value2 = 8
$ws.Range("D6").Value2 = 16
$ws.Range("D7").Value2 = 32

# Bottom border row D8 takes the same style as C8.
$ws.Range("C8").Copy()
$ws.Range("D8").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# Step 4: fill in the new Uint32 row contents (row 7), which currently
# only has the Bit Size value set above.
# ---------------------------------------------------------------------
$ws.Range("B6").Copy()
$ws.Range("B7").PasteSpecial($xlPasteFormats)
$ws.Range("B7").Value2 = "Uint32"

$ws.Range("C6").Copy()
$ws.Range("C7").PasteSpecial($xlPasteFormats)
$ws.Range("C7").Value2 = "Base Types"

$ws.Range("E6").Copy()
$ws.Range("E7").PasteSpecial($xlPasteFormats)
$ws.Range("E7").Value2 = "unsigned long"
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# Step 5: column width - the Bit Size column is narrower than the old
# Native Declaration column used to be.
# ---------------------------------------------------------------------
$ws.Columns("D").ColumnWidth = 21.875

# ---------------------------------------------------------------------
# Step 6: restore the view - active selection and window position.
# ---------------------------------------------------------------------
$ws.Range("I17").Select()
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 1

Write-Host "done"
